$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("D2").Value = 13
$ws1.Range("H2").Value = 0.08
$ws1.Range("L2").Value = 1.07

$ws1.Range("D3").Value = 13
$ws1.Range("L3").Value = 1.17

$ws1.Range("L4").Value = 0.93
$ws1.Range("L5").Value = 0.96
$ws1.Range("L6").Value = 1.14
$ws1.Range("L7").Value = 1.08
$ws1.Range("L8").Value = 0.87
$ws1.Range("L9").Value = 1.06
$ws1.Range("L10").Value = 1.16
$ws1.Range("L11").Value = 0.8
$ws1.Range("L12").Value = 0.87
$ws1.Range("L13").Value = 0.89
$ws1.Range("L14").Value = 0.89
$ws1.Range("L15").Value = 1.03
$ws1.Range("L16").Value = 0.9
$ws1.Range("L17").Value = 0.86

# --- Sheet: Summary ---
# B9:B12 hold numeric-looking values but must remain stored as TEXT
# (matching the source inlineStr type). Directly assigning a numeric
# string makes Excel auto-coerce it to a number, so we briefly force
# a Text number format, assign the value, then strip the formatting
# back off so the cell keeps its original (default) style.
$ws2 = $wb.Worksheets.Item("Summary")

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws2.Range("B9") "174"
Set-TextValue $ws2.Range("B10") "94"
Set-TextValue $ws2.Range("B11") "51"
Set-TextValue $ws2.Range("B12") "13"
